$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data (row 3), A3:G3 = "eddie"
$ws.Range("A3:G3").Value = "eddie"

# Update the selection to match the diff (activeCell F8, sqref F8)
$ws.Range("F8").Select()
